$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title heading and matching bold line near the end (two occurrences, same replacement)
Replace-Text "Play Knockout Football Rush for free - review and gameplay" "Play Knockout Football Rush Free - Review"

# "What we like" bullet list
Replace-Text "Multiple ways to win both horizontally and vertically" "Unique gameplay mechanics with both vertical and horizontal wins"
Replace-Text "Exciting penalty kick bonus feature with rewards up to 60x stake" "Modern and visually appealing graphic design"
Replace-Text "Modern design and sound effects that enhance the soccer experience" "Multiple betting options for players"
Replace-Text "Flexible betting options for different budgets" "Exciting bonus features with potential for big rewards"

# "What we don't like" bullet list
Replace-Text "High volatility may not suit all players" "Limited availability of the joker character on the 3rd row"
Replace-Text "Limited bonus features besides the penalty kick feature" "High volatility may not appeal to all players"

# Closing italic summary paragraph
Replace-Text "Read our review of the soccer-themed slot Knockout Football Rush. Play for free and experience exciting bonus features like penalty kicks and wild multipliers." "Read our review of Knockout Football Rush and play this exciting slot game for free."
